{"js": "// \"Correciones trazos finos de ventas\"\n//\n// Four sentences in the use-case table are reworded:\n//   1. \"El EV selecciona cada producto que desea agregar al cat\u00e1logo\"\n//        -> \"El EV selecciona los producto que desea agregar al cat\u00e1logo\"\n//   2. \"El sistema solicita si se desea confirmar el cat\u00e1logo\"\n//        -> \"El sistema solicita si se desea confirmar la creaci\u00f3n del cat\u00e1logo\"\n//   3. \"El EV confirma el cat\u00e1logo\"\n//        -> \"El EV confirma la creaci\u00f3n del cat\u00e1logo.\"\n//   4. \"El EV no confirma el cat\u00e1logo\"\n//        -> \"El EV no confirma la creaci\u00f3n del cat\u00e1logo.\"\n//\n// Each original sentence lives in its own unique run, so we can find it with\n// Body.search() (exact, case-sensitive) and rewrite its text in place.\n\nconst body = context.document.body;\n\nconst edits = [\n  {\n    find: \"El EV selecciona cada producto que desea agregar al cat\u00e1logo\",\n    replace: \"El EV selecciona los producto que desea agregar al cat\u00e1logo\"\n  },\n  {\n    find: \"El sistema solicita si se desea confirmar el cat\u00e1logo\",\n    replace: \"El sistema solicita si se desea confirmar la creaci\u00f3n del cat\u00e1logo\"\n  },\n  {\n    find: \"El EV confirma el cat\u00e1logo\",\n    replace: \"El EV confirma la creaci\u00f3n del cat\u00e1logo.\"\n  },\n  {\n    find: \"El EV no confirma el cat\u00e1logo\",\n    replace: \"El EV no confirma la creaci\u00f3n del cat\u00e1logo.\"\n  }\n];\n\nfor (const { find, replace } of edits) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# \"Correciones trazos finos de ventas\"\n#\n# Four sentences in the use-case table are reworded:\n#   1. \"El EV selecciona cada producto que desea agregar al cat\u00e1logo\"\n#        -> \"El EV selecciona los producto que desea agregar al cat\u00e1logo\"\n#   2. \"El sistema solicita si se desea confirmar el cat\u00e1logo\"\n#        -> \"El sistema solicita si se desea confirmar la creaci\u00f3n del cat\u00e1logo\"\n#   3. \"El EV confirma el cat\u00e1logo\"\n#        -> \"El EV confirma la creaci\u00f3n del cat\u00e1logo.\"\n#   4. \"El EV no confirma el cat\u00e1logo\"\n#        -> \"El EV no confirma la creaci\u00f3n del cat\u00e1logo.\"\n#\n# Word constants used below (WdFindWrap.wdFindContinue / WdReplace.wdReplaceOne)\n$wdFindContinue = 1\n$wdReplaceOne = 2\n\n$d = $word.ActiveDocument\n\n$edits = @(\n    @{ Find = \"El EV selecciona cada producto que desea agregar al cat\u00e1logo\"; Replace = \"El EV selecciona los producto que desea agregar al cat\u00e1logo\" },\n    @{ Find = \"El sistema solicita si se desea confirmar el cat\u00e1logo\";        Replace = \"El sistema solicita si se desea confirmar la creaci\u00f3n del cat\u00e1logo\" },\n    @{ Find = \"El EV confirma el cat\u00e1logo\";                                  Replace = \"El EV confirma la creaci\u00f3n del cat\u00e1logo.\" },\n    @{ Find = \"El EV no confirma el cat\u00e1logo\";                               Replace = \"El EV no confirma la creaci\u00f3n del cat\u00e1logo.\" }\n)\n\nforeach ($edit in $edits) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute(\n        $edit.Find,    # FindText\n        $true,         # MatchCase\n        $true,         # MatchWholeWord\n        $false,        # MatchWildcards\n        $false,        # MatchSoundsLike\n        $false,        # MatchAllWordForms\n        $true,         # Forward\n        $wdFindContinue, # Wrap\n        $false,        # Format\n        $edit.Replace, # ReplaceWith\n        $wdReplaceOne  # Replace\n    )\n\n    if (-not $found) {\n        throw \"Text not found: $($edit.Find)\"\n    }\n}\n"}
